# Auto-generated Excel COM-interop script to apply Shinryu_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2571.4614
$ws.Range("I132").Value = 2646.652
$ws.Range("J132").Value = 1995
$ws.Range("K132").Value = 7939.956
$ws.Range("L132").Value = 5985
$ws.Range("M132").Value = -5409.956
$ws.Range("N132").Value = -11045

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 7377
$ws.Range("J43").Value = 7377
$ws.Range("L43").Value = 7377
$ws.Range("N43").Value = -8003
$ws.Range("H45").Value = 2126.111
$ws.Range("I45").Value = 2141.875
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2141.875
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1764.875
$ws.Range("N45").Value = -2754
$ws.Range("H74").Value = 7297.9414
$ws.Range("I74").Value = 7670.067
$ws.Range("J74").Value = 4507
$ws.Range("K74").Value = 7670.067
$ws.Range("L74").Value = 4507
$ws.Range("M74").Value = -6796.067
$ws.Range("N74").Value = -6255
$ws.Range("H77").Value = 7297.9414
$ws.Range("I77").Value = 7670.067
$ws.Range("J77").Value = 4507
$ws.Range("K77").Value = 38350.335
$ws.Range("L77").Value = 22535
$ws.Range("M77").Value = -33982.335
$ws.Range("N77").Value = -31271
$ws.Range("H109").Value = 28399.6
$ws.Range("J109").Value = 28399.6
$ws.Range("L109").Value = 28399.6
$ws.Range("N109").Value = -31173.6
$ws.Range("H110").Value = 1029.6285
$ws.Range("I110").Value = 742.96295
$ws.Range("J110").Value = 1997.125
$ws.Range("K110").Value = 742.96295
$ws.Range("L110").Value = 1997.125
$ws.Range("M110").Value = 1302.03705
$ws.Range("N110").Value = -6087.125
$ws.Range("H122").Value = 2268.5715
$ws.Range("I122").Value = 1976
$ws.Range("K122").Value = 5928
$ws.Range("M122").Value = -3478
$ws.Range("H132").Value = 1880.5264
$ws.Range("I132").Value = 1257.5714
$ws.Range("J132").Value = 3624.8
$ws.Range("K132").Value = 3772.7142
$ws.Range("L132").Value = 10874.4
$ws.Range("M132").Value = -1242.7142
$ws.Range("N132").Value = -15934.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 172.5
$ws.Range("I80").Value = 203
$ws.Range("J80").Value = 121.666664
$ws.Range("K80").Value = 203
$ws.Range("L80").Value = 121.666664
$ws.Range("M80").Value = 795
$ws.Range("N80").Value = -2117.666664
$ws.Range("H83").Value = 172.5
$ws.Range("I83").Value = 203
$ws.Range("J83").Value = 121.666664
$ws.Range("K83").Value = 1015
$ws.Range("L83").Value = 608.33332
$ws.Range("M83").Value = 3977
$ws.Range("N83").Value = -10592.33332
$ws.Range("H99").Value = 1844.55
$ws.Range("I99").Value = 1521.625
$ws.Range("J99").Value = 3136.25
$ws.Range("K99").Value = 1521.625
$ws.Range("L99").Value = 3136.25
$ws.Range("M99").Value = -23.625
$ws.Range("N99").Value = -6132.25
$ws.Range("H105").Value = 2776.61
$ws.Range("I105").Value = 1802.375
$ws.Range("J105").Value = 2861.3262
$ws.Range("K105").Value = 1802.375
$ws.Range("L105").Value = 2861.3262
$ws.Range("M105").Value = -55.375
$ws.Range("N105").Value = -6355.3262
$ws.Range("H134").Value = 1843.9131
$ws.Range("I134").Value = 1684
$ws.Range("J134").Value = 2603.5
$ws.Range("K134").Value = 5052
$ws.Range("L134").Value = 7810.5
$ws.Range("M134").Value = -2517
$ws.Range("N134").Value = -12880.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 525.1177
$ws.Range("I107").Value = 463.39285
$ws.Range("J107").Value = 813.1667
$ws.Range("K107").Value = 463.39285
$ws.Range("L107").Value = 813.1667
$ws.Range("M107").Value = 1456.60715
$ws.Range("N107").Value = -4653.1667
$ws.Range("H122").Value = 948.619
$ws.Range("I122").Value = 896.05
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2688.15
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -238.1499999999996
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 4509.75
$ws.Range("I132").Value = 3175.2
$ws.Range("J132").Value = 5463
$ws.Range("K132").Value = 9525.599999999999
$ws.Range("L132").Value = 16389
$ws.Range("M132").Value = -6995.599999999999
$ws.Range("N132").Value = -21449
$ws.Range("H134").Value = 2681.6667
$ws.Range("I134").Value = 1462.2941
$ws.Range("J134").Value = 7864
$ws.Range("K134").Value = 4386.8823
$ws.Range("L134").Value = 23592
$ws.Range("M134").Value = -1851.8823
$ws.Range("N134").Value = -28662

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 5006.8
$ws.Range("I101").Value = 5034
$ws.Range("K101").Value = 15102
$ws.Range("M101").Value = -12668
$ws.Range("H131").Value = 857.7931
$ws.Range("I131").Value = 544.6667
$ws.Range("J131").Value = 915.30615
$ws.Range("K131").Value = 1634.0001
$ws.Range("L131").Value = 2745.91845
$ws.Range("M131").Value = 3405.9999
$ws.Range("N131").Value = -12825.91845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2033.1
$ws.Range("I102").Value = 1981.4736
$ws.Range("J102").Value = 3014
$ws.Range("K102").Value = 1981.4736
$ws.Range("L102").Value = 3014
$ws.Range("M102").Value = -359.4736
$ws.Range("N102").Value = -6258
$ws.Range("H126").Value = 3693.7144
$ws.Range("I126").Value = 3406.96
$ws.Range("J126").Value = 4410.6
$ws.Range("K126").Value = 10220.88
$ws.Range("L126").Value = 13231.8
$ws.Range("M126").Value = -7750.880000000001
$ws.Range("N126").Value = -18171.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2644.0557
$ws.Range("I7").Value = 2379.4666
$ws.Range("J7").Value = 3967
$ws.Range("K7").Value = 2379.4666
$ws.Range("L7").Value = 3967
$ws.Range("M7").Value = -2267.4666
$ws.Range("N7").Value = -4191
$ws.Range("H126").Value = 2644.0557
$ws.Range("I126").Value = 2379.4666
$ws.Range("J126").Value = 3967
$ws.Range("K126").Value = 7138.399800000001
$ws.Range("L126").Value = 11901
$ws.Range("M126").Value = -4668.399800000001
$ws.Range("N126").Value = -16841

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2682.077
$ws.Range("I122").Value = 2488.9167
$ws.Range("K122").Value = 7466.750100000001
$ws.Range("M122").Value = -5016.750100000001
